# Renumber [[PERSON_N]] placeholders across the two numbered lists so that
# repeated references within a paragraph collapse onto a single canonical
# identity (fixing the -ina/-ynou declension fallback corruption).
$d = $word.ActiveDocument

$pairs = @(
    ,@("[[PERSON_3]] – „pro [[PERSON_4]]“, „s [[PERSON_5]]“", "[[PERSON_3]] – „pro [[PERSON_3]]“, „s [[PERSON_3]]“")
    ,@("[[PERSON_6]] – „s [[PERSON_6]]“, „o [[PERSON_6]]“", "[[PERSON_4]] – „s [[PERSON_4]]“, „o [[PERSON_4]]“")
    ,@("[[PERSON_7]] – „u [[PERSON_7]]“, „k [[PERSON_7]]“", "[[PERSON_5]] – „u [[PERSON_5]]“, „k [[PERSON_5]]“")
    ,@("[[PERSON_8]] – „o [[PERSON_8]]“, „se [[PERSON_8]]“", "[[PERSON_6]] – „o [[PERSON_6]]“, „se [[PERSON_6]]“")
    ,@("[[PERSON_9]] – „k [[PERSON_9]]“, „u [[PERSON_9]]“", "[[PERSON_7]] – „k [[PERSON_7]]“, „u [[PERSON_7]]“")
    ,@("[[PERSON_10]] – „s [[PERSON_10]]“, „o [[PERSON_10]]“", "[[PERSON_8]] – „s [[PERSON_8]]“, „o [[PERSON_8]]“")
    ,@("[[PERSON_11]] – „u [[PERSON_11]]“, „s [[PERSON_11]]“", "[[PERSON_9]] – „u [[PERSON_9]]“, „s [[PERSON_9]]“")
    ,@("[[PERSON_12]] – „s [[PERSON_12]]“, „k [[PERSON_12]]“", "[[PERSON_10]] – „s [[PERSON_10]]“, „k [[PERSON_10]]“")
    ,@("[[PERSON_13]] – „s [[PERSON_14]]“, „o [[PERSON_15]]“", "[[PERSON_11]] – „s [[PERSON_12]]“, „o [[PERSON_13]]“")
    ,@("[[PERSON_16]] – „ke [[PERSON_16]]“, „o [[PERSON_16]]“", "[[PERSON_14]] – „ke [[PERSON_14]]“, „o [[PERSON_14]]“")
    ,@("[[PERSON_17]] – „o [[PERSON_17]]“, „s [[PERSON_17]]“", "[[PERSON_15]] – „o [[PERSON_15]]“, „s [[PERSON_15]]“")
    ,@("[[PERSON_18]] – „u [[PERSON_18]]“, „s [[PERSON_18]]“", "[[PERSON_16]] – „u [[PERSON_16]]“, „s [[PERSON_16]]“")
    ,@("[[PERSON_19]] – „ke [[PERSON_19]]“, „o [[PERSON_19]]“", "[[PERSON_17]] – „ke [[PERSON_17]]“, „o [[PERSON_17]]“")
    ,@("[[PERSON_20]] – „s [[PERSON_21]]“, „o [[PERSON_20]]“", "[[PERSON_18]] – „s [[PERSON_19]]“, „o [[PERSON_18]]“")
    ,@("[[PERSON_22]] – „s [[PERSON_22]]“, „o [[PERSON_22]]“", "[[PERSON_20]] – „s [[PERSON_20]]“, „o [[PERSON_20]]“")
    ,@("[[PERSON_23]] – „k [[PERSON_23]]“, „od [[PERSON_23]]“", "[[PERSON_21]] – „k [[PERSON_21]]“, „od [[PERSON_21]]“")
    ,@("[[PERSON_24]] – „o [[PERSON_24]]“, „s [[PERSON_24]]“", "[[PERSON_22]] – „o [[PERSON_22]]“, „s [[PERSON_22]]“")
    ,@("[[PERSON_25]] – „o [[PERSON_26]]“, „se [[PERSON_27]]“", "[[PERSON_23]] – „o [[PERSON_24]]“, „se [[PERSON_25]]“")
    ,@("[[PERSON_28]] – „s [[PERSON_28]]“, „u [[PERSON_28]]“", "[[PERSON_26]] – „s [[PERSON_26]]“, „u [[PERSON_26]]“")
    ,@("[[PERSON_29]] – „o [[PERSON_30]]“, „s [[PERSON_31]]“", "[[PERSON_27]] – „o [[PERSON_28]]“, „s [[PERSON_29]]“")
    ,@("[[PERSON_32]] – „k [[PERSON_32]]“, „o [[PERSON_32]]“", "[[PERSON_30]] – „k [[PERSON_30]]“, „o [[PERSON_30]]“")
    ,@("[[PERSON_33]] – „se [[PERSON_33]]“, „o Soně Mikulkové“", "[[PERSON_31]] – „se [[PERSON_31]]“, „o Soně Mikulkové“")
    ,@("[[PERSON_34]] – „o [[PERSON_34]]“, „s [[PERSON_34]]“", "[[PERSON_32]] – „o [[PERSON_32]]“, „s [[PERSON_32]]“")
    ,@("[[PERSON_35]] – „s [[PERSON_36]]“, „o [[PERSON_35]]“", "[[PERSON_33]] – „s [[PERSON_34]]“, „o [[PERSON_33]]“")
    ,@("[[PERSON_37]] – „k [[PERSON_38]]“, „s [[PERSON_37]]“", "[[PERSON_35]] – „k [[PERSON_36]]“, „s [[PERSON_35]]“")
    ,@("[[PERSON_39]] – „s [[PERSON_39]]“, „o [[PERSON_40]]“", "[[PERSON_37]] – „s [[PERSON_37]]“, „o [[PERSON_38]]“")
    ,@("[[PERSON_41]] – „od [[PERSON_41]]“, „s [[PERSON_41]]“", "[[PERSON_39]] – „od [[PERSON_39]]“, „s [[PERSON_39]]“")
    ,@("[[PERSON_42]] – „k [[PERSON_43]]“, „o [[PERSON_43]]“", "[[PERSON_40]] – „k [[PERSON_41]]“, „o [[PERSON_41]]“")
    ,@("[[PERSON_44]] – „o [[PERSON_45]]“, „s [[PERSON_44]]“", "[[PERSON_42]] – „o [[PERSON_43]]“, „s [[PERSON_42]]“")
    ,@("[[PERSON_48]] – „s [[PERSON_48]]“, „o [[PERSON_49]]“", "[[PERSON_44]] – „s [[PERSON_44]]“, „o [[PERSON_45]]“")
    ,@("[[PERSON_50]] – „k [[PERSON_50]]“, „s [[PERSON_51]]“", "[[PERSON_48]] – „k [[PERSON_48]]“, „s [[PERSON_49]]“")
    ,@("[[PERSON_52]] – „pro [[PERSON_53]]“, „o [[PERSON_54]]“", "[[PERSON_50]] – „pro [[PERSON_51]]“, „o [[PERSON_52]]“")
    ,@("[[PERSON_55]] – „k [[PERSON_55]]“, „o [[PERSON_55]]“", "[[PERSON_53]] – „k [[PERSON_53]]“, „o [[PERSON_53]]“")
    ,@("[[PERSON_56]] – „o [[PERSON_57]]“, „s [[PERSON_56]]“", "[[PERSON_54]] – „o [[PERSON_55]]“, „s [[PERSON_54]]“")
    ,@("[[PERSON_58]] – „s [[PERSON_59]]“, „o [[PERSON_60]]“", "[[PERSON_56]] – „s [[PERSON_57]]“, „o [[PERSON_58]]“")
    ,@("[[PERSON_61]] – „s [[PERSON_61]]“, „o [[PERSON_61]]“", "[[PERSON_59]] – „s [[PERSON_59]]“, „o [[PERSON_59]]“")
    ,@("[[PERSON_62]] – „u [[PERSON_62]]“, „o [[PERSON_63]]“", "[[PERSON_60]] – „u [[PERSON_60]]“, „o [[PERSON_61]]“")
    ,@("[[PERSON_64]] – „se [[PERSON_65]]“, „o [[PERSON_64]]“", "[[PERSON_62]] – „se [[PERSON_63]]“, „o [[PERSON_62]]“")
    ,@("[[PERSON_66]] – „o [[PERSON_67]]“, „s [[PERSON_68]]“", "[[PERSON_64]] – „o [[PERSON_65]]“, „s [[PERSON_66]]“")
    ,@("[[PERSON_69]] – „k [[PERSON_70]]“, „o [[PERSON_70]]“", "[[PERSON_67]] – „k [[PERSON_68]]“, „o [[PERSON_68]]“")
    ,@("[[PERSON_71]] – „o [[PERSON_72]]“, „s [[PERSON_71]]“", "[[PERSON_69]] – „o [[PERSON_70]]“, „s [[PERSON_69]]“")
    ,@("[[PERSON_73]] – „s [[PERSON_73]]“, „o [[PERSON_73]]“", "[[PERSON_71]] – „s [[PERSON_71]]“, „o [[PERSON_71]]“")
    ,@("[[PERSON_74]] – „s [[PERSON_74]]“, „o [[PERSON_75]]“", "[[PERSON_72]] – „s [[PERSON_72]]“, „o [[PERSON_73]]“")
    ,@("[[PERSON_76]] – „o [[PERSON_77]]“, „s [[PERSON_78]]“", "[[PERSON_74]] – „o [[PERSON_75]]“, „s [[PERSON_76]]“")
    ,@("[[PERSON_79]] – „s [[PERSON_79]]“, „o [[PERSON_80]]“", "[[PERSON_77]] – „s [[PERSON_77]]“, „o [[PERSON_78]]“")
    ,@("[[PERSON_81]] – „o [[PERSON_81]]“, „s [[PERSON_81]]“", "[[PERSON_79]] – „o [[PERSON_80]]“, „s [[PERSON_79]]“")
    ,@("[[PERSON_82]] – „s [[PERSON_82]]“, „o [[PERSON_83]]“", "[[PERSON_81]] – „s [[PERSON_81]]“, „o [[PERSON_82]]“")
)

$replacements = 0
$notFound = @()
$i = 0
foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if ($found) { $replacements++ } else { $notFound += $i }
    $i++
}

Write-Output ("Replacements applied: " + $replacements + " / " + $pairs.Count)
if ($notFound.Count -gt 0) {
    Write-Output ("WARNING - not found indices: " + ($notFound -join ","))
}
